$d = $word.ActiveDocument

function Replace-InParagraph {
    param($Index, $OldText, $NewText)
    $p = $d.Paragraphs.Item($Index)
    $rng = $p.Range
    $ok = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 0, $false, $NewText, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed for paragraph" $Index
    }
}

Replace-InParagraph 19 "1. Rocha Filho, J.A., Vitolo, M. Guia para aulas práticas de biotecnologia de enzimas e fermentação. Editora Blucher, 2021. `v2. Cisternas, J.R. Fundamentos de bioquímica experimental. São Paulo: Atheneu, 2005. `v3. Nelson, D.L., Cox, M.M. Princípios de bioquímica de Lehninger. Artmed Editora, 2022. `v4. Voet, D., Voet, J.G., Pratt, C.W. Fundamentos de Bioquímica: a vida em nivel molecular. Artmed Editora, 2014. `v5. Vitolo, M., Pessoa Junior, A., Monteiro, G., Carvalho, J.C.M., Stephano, M.A., Sato, S. Biotecnologia farmacêutica: aspectos sobre aplicação industrial. Editora Blucher, 2015." "6007846 - Júlio César dos Santos"
Replace-InParagraph 17 "A recuperação será realizada através de uma prova escrita (PR) e a média de recuperação (MR) será calculada conforme: MR = (NF + PR)/2." "5082401 - André Moreni Lopes"
Replace-InParagraph 17 "A nota final (NF) será calculada conforme: NF = (P + R)/2. A" "1. Rocha Filho, J.A., Vitolo, M. Guia para aulas práticas de biotecnologia de enzimas e fermentação. Editora Blucher, 2021. `v2. Cisternas, J.R. Fundamentos de bioquímica experimental. São Paulo: Atheneu, 2005. `v3. Nelson, D.L., Cox, M.M. Princípios de bioquímica de Lehninger. Artmed Editora, 2022. `v4. Voet, D., Voet, J.G., Pratt, C.W. Fundamentos de Bioquímica: a vida em nivel molecular. Artmed Editora, 2014. `v5. Vitolo, M., Pessoa Junior, A., Monteiro, G., Carvalho, J.C.M., Stephano, M.A., Sato, S. Biotecnologia farmacêutica: aspectos sobre aplicação industrial. Editora Blucher, 2015."
Replace-InParagraph 17 "A avaliação será realizada através de uma prova escrita (P) e um relatório de atividades práticas (R)." "A recuperação será realizada através de uma prova escrita (PR) e a média de recuperação (MR) será calculada conforme: MR = (NF + PR)/2."
Replace-InParagraph 6 "Desenvolver e aperfeiçoar o entendimento teórico e prático dos processos bioquímicos fundamentais através da realização de atividades práticas de laboratório." "Reação de saponificação; Enzimas proteolíticas em produtos comerciais; Extração líquido-líquido de proteínas; Biomateriais sustentáveis; Produção e destilação de etanol; Precipitação de biomoléculas e Reação de Hill."
Replace-InParagraph 7 "Developing and enhancing the theoretical and practical understanding of fundamental biochemical processes through the execution of laboratory practical activities." "Saponification reaction; Proteolytic enzymes in commercial products; Liquid-liquid extraction of proteins; Sustainable biomaterials; Production and distillation of ethanol; Precipitation of biomolecules; and Hill reaction."
Replace-InParagraph 9 "5082401 - André Moreni Lopes`v6007846 - Júlio César dos Santos" "Desenvolver e aperfeiçoar o entendimento teórico e prático dos processos bioquímicos fundamentais através da realização de atividades práticas de laboratório.`vAplicação da saponificação em processos industriais, agentes envolvidos na reação e sua utilização em produtos comerciais (cálculo de rendimento, CMC e pH). Avaliação enzimática de proteases para determinação de sua atividade proteolítica em produtos comerciais (sabão em pó, detergentes e cosméticos). Extração líquido-líquido de proteínas e enzimas utilizando solventes orgânicos e polímeros/tensoativos - quantificação dos parâmetros de extração (balanço de massa, recuperação, fator de purificação). Obtenção de biomateriais (bioplástico) de interesse biotecnológico derivado de fontes biológicas - cálculo do rendimento; caracterização do produto final obtido (textura, cor e cheiro) e comparação com os plásticos convencionais. Produção e destilação de etanol - conceitos gerais e fermentação de glicose; produção de etanol e CO2; consumo da fonte de carbono; cálculo da eficiência do processo; ação de um inibidor da glicólise. Precipitação de biomoléculas utilizando diferentes agentes precipitadores (sais, polímeros e solventes orgânicos) - quantificação da recuperação, pH e potencial elétrico. Extração de clorofila e reação de Hill - estrutura de cloroplastos; papel da clorofila nos sistemas fotossintéticos I e II; fase escura/luminosa; produção de NADP; produção de ATP; papel do corante como aceptor de prótons e elétrons. *Dentro do programa da disciplina é planejado realizar eventual `"Visita Didática Complementar`"."
Replace-InParagraph 11 "Reação de saponificação; Enzimas proteolíticas em produtos comerciais; Extração líquido-líquido de proteínas; Biomateriais sustentáveis; Produção e destilação de etanol; Precipitação de biomoléculas e Reação de Hill." "A avaliação será realizada através de uma prova escrita (P) e um relatório de atividades práticas (R)."
Replace-InParagraph 12 "Saponification reaction; Proteolytic enzymes in commercial products; Liquid-liquid extraction of proteins; Sustainable biomaterials; Production and distillation of ethanol; Precipitation of biomolecules; and Hill reaction." "Developing and enhancing the theoretical and practical understanding of fundamental biochemical processes through the execution of laboratory practical activities."
Replace-InParagraph 14 "Aplicação da saponificação em processos industriais, agentes envolvidos na reação e sua utilização em produtos comerciais (cálculo de rendimento, CMC e pH). Avaliação enzimática de proteases para determinação de sua atividade proteolítica em produtos comerciais (sabão em pó, detergentes e cosméticos). Extração líquido-líquido de proteínas e enzimas utilizando solventes orgânicos e polímeros/tensoativos - quantificação dos parâmetros de extração (balanço de massa, recuperação, fator de purificação). Obtenção de biomateriais (bioplástico) de interesse biotecnológico derivado de fontes biológicas - cálculo do rendimento; caracterização do produto final obtido (textura, cor e cheiro) e comparação com os plásticos convencionais. Produção e destilação de etanol - conceitos gerais e fermentação de glicose; produção de etanol e CO2; consumo da fonte de carbono; cálculo da eficiência do processo; ação de um inibidor da glicólise. Precipitação de biomoléculas utilizando diferentes agentes precipitadores (sais, polímeros e solventes orgânicos) - quantificação da recuperação, pH e potencial elétrico. Extração de clorofila e reação de Hill - estrutura de cloroplastos; papel da clorofila nos sistemas fotossintéticos I e II; fase escura/luminosa; produção de NADP; produção de ATP; papel do corante como aceptor de prótons e elétrons. *Dentro do programa da disciplina é planejado realizar eventual `"Visita Didática Complementar`"." "A nota final (NF) será calculada conforme: NF = (P + R)/2. A"

Write-Host "All replacements applied."
